$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "10 Week"
$ws2 = $wb.Worksheets.Item(2)   # "Daily report"

# ---------------------------------------------------------------------------
# Sheet2 ("Daily report"): new daily-log entry for row 6 (week 5 audio work)
# ---------------------------------------------------------------------------
$ws2.Range("C6").Value = "Generated all audio. Organized RIR analysis folders. Modified MATLAB scripts for compatibility. Investigated and troubleshooted EDT/RT60 analysis issues. Experimented with different audio settings (sampling rate, sine sweep type, probe number, audio delay, volume). Identified persistent problems with audio analysis results."
$ws2.Range("A6:G6").RowHeight = 101.5

# ---------------------------------------------------------------------------
# Sheet1 ("10 Week"): Gantt-chart status updates + extended legend
# ---------------------------------------------------------------------------

# Weekly note for week 5 + new daily-log note cross references
$ws1.Range("N5").Value = "Gonna go longer than week 4, hopefully finish analysis with pretty graphs etc by end of Week 5"

# New legend entries: "Extended" (orange/theme) swatch and "Cancelled/Delayed" (red) swatch
$ws1.Range("R1").Value = "Extended"
$ws1.Range("T1").Value = "Cancelled/Delayed"

$ws1.Range("N6").Value = "This is now optional, SSC Mona is using stereo (some problem/setback?)"

# Mark previously "Todo" (blue) cells as "Done" (green)
$ws1.Range("Q1").Copy()
$ws1.Range("B3").PasteSpecial(-4122)
$ws1.Range("C4").PasteSpecial(-4122)
$ws1.Range("D4").PasteSpecial(-4122)
$ws1.Range("D5").PasteSpecial(-4122)
$ws1.Range("E5").PasteSpecial(-4122)

# Build the "Extended" fill cleanly (reuse an existing solid fill as the base so
# the engine doesn't leave a transient blank fill behind), then recolor to theme accent2.
$ws1.Range("Q1").Copy()
$ws1.Range("F5").PasteSpecial(-4122)
$ws1.Range("F5").Interior.ThemeColor = 6

# Build the "Cancelled/Delayed" (red) fill cleanly the same way.
$ws1.Range("Q1").Copy()
$ws1.Range("F6").PasteSpecial(-4122)
$ws1.Range("F6").Interior.Color = 255

# Reuse those two freshly-created styles for the legend swatches.
$ws1.Range("F5").Copy()
$ws1.Range("S1").PasteSpecial(-4122)
$ws1.Range("F6").Copy()
$ws1.Range("V1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selections / active sheet: "Daily report" selection moves to C5, then "10 Week"
# becomes the active tab again with N6 selected.
# ---------------------------------------------------------------------------
$ws2.Range("C5").Select()
$ws1.Activate()
$ws1.Range("N6").Select()
